$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each data row (2..37) in this "Pepino dulce" weekly price sheet, the
# observation's date/quality/volume/price fields (columns D, I, J, K, L, M, P)
# get reshuffled among the rows (a weekly re-ordering), while the rest of the
# row (market, region, category, unit, origin, etc.) stays in place.
#
# $mapping[newRow] = oldRow  -- the row whose D/I/J/K/L/M/P values should be
# placed into newRow.
$mapping = @{ 2=3; 3=37; 4=24; 5=33; 6=12; 7=4; 8=9; 9=10; 10=23; 11=32; 12=26; 13=16; 14=18; 15=36; 16=19; 17=2; 18=34; 19=27; 20=7; 21=29; 22=30; 23=6; 24=11; 25=28; 26=21; 27=20; 28=17; 29=14; 30=15; 31=22; 32=13; 33=31; 34=8; 35=25; 36=35; 37=5 }

# Capture the original (pre-edit) values for columns D, I, J, K, L, M, P for
# every data row before any writes happen, so the reshuffle reads are not
# affected by earlier writes.
$origD = @{}
$origI = @{}
$origJ = @{}
$origK = @{}
$origL = @{}
$origM = @{}
$origP = @{}

for ($r = 2; $r -le 37; $r++) {
    $origD[$r] = $ws.Cells.Item($r, 4).Value()
    $origI[$r] = $ws.Cells.Item($r, 9).Value()
    $origJ[$r] = $ws.Cells.Item($r, 10).Value()
    $origK[$r] = $ws.Cells.Item($r, 11).Value()
    $origL[$r] = $ws.Cells.Item($r, 12).Value()
    $origM[$r] = $ws.Cells.Item($r, 13).Value()
    $origP[$r] = $ws.Cells.Item($r, 16).Value()
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $ws.Cells.Item($newRow, 4).Value = $origD[$oldRow]
    $ws.Cells.Item($newRow, 9).Value = $origI[$oldRow]
    $ws.Cells.Item($newRow, 10).Value = $origJ[$oldRow]
    $ws.Cells.Item($newRow, 11).Value = $origK[$oldRow]
    $ws.Cells.Item($newRow, 12).Value = $origL[$oldRow]
    $ws.Cells.Item($newRow, 13).Value = $origM[$oldRow]
    $ws.Cells.Item($newRow, 16).Value = $origP[$oldRow]
}
